# FY2026-自治体生成AI案件-フェーズ別見通し.pptx
# Slide 3 ("落札金額" phase-breakdown table) update:
#   1. Widen + re-word the section title textbox.
#   2. Append a second line "(判明件数/総件数件)" under every award-amount
#      figure in the breakdown table, shrinking the font from 10pt to 9pt
#      so the two-line cells still fit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# ---- 1. Title textbox: widen and reword -------------------------------
$title = $s.Shapes.Item(4)
$title.Width = 432
$title.TextFrame.TextRange.Text = "落札金額 合計（カッコ内は金額判明件数）"

# ---- 2. Data table: add a second "(判明/総数件)" line to each figure --
$tbl = $s.Shapes.Item(5).Table

# counts[row][col] -> text to append as a new paragraph, keyed by the
# table's 1-based Row/Column indices (row 1 is the header, column 1 is
# the fiscal-year label - neither is touched).
$counts = @{
    "2,2" = "(1/4件)";    "2,3" = "(0/6件)";   "2,4" = "(5/6件)";
    "2,5" = "(40/72件)";  "2,6" = "(25/42件)"; "2,7" = "(71/130件)";

    "3,2" = "(7/15件)";   "3,3" = "(0/6件)";   "3,4" = "(7/12件)";
    "3,5" = "(60/132件)"; "3,6" = "(43/67件)"; "3,7" = "(117/232件)";

    "4,2" = "(10/18件)";  "4,3" = "(0/8件)";   "4,4" = "(4/7件)";
    "4,5" = "(44/130件)"; "4,6" = "(42/62件)"; "4,7" = "(100/225件)";
}

for ($r = 2; $r -le 4; $r++) {
    for ($c = 2; $c -le 7; $c++) {
        $key = "$r,$c"
        $suffix = $counts[$key]
        $cell = $tbl.Cell($r, $c)
        $tr = $cell.Shape.TextFrame.TextRange
        $tr.Text = $tr.Text + "`r" + $suffix
        $tr.Font.Size = 9
    }
}
